$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Rights")

# Insert a new row at position 110 ("DOCUMENT_TEMPLATE_MANAGEMENT"), pushing
# the existing rows 110-127 down to 111-128.
$ws.Rows.Item(110).Insert()

# The freshly inserted row inherits formatting from the row above (row 109),
# which is itself an all-"No" row with a bold A/B label style - fine for A/B,
# but we need the correct Yes/No fill styles on C:Y. Pull those styles from a
# known-good Yes cell (C109) and No cell (D109) elsewhere in the sheet so the
# new row reuses the same style indices as the rest of the table.
$ws.Range("C109").Copy()
$ws.Range("C110").PasteSpecial(-4122)

$ws.Range("D109").Copy()
$ws.Range("D110:Y110").PasteSpecial(-4122)

$ws.Cells.Item(110, 1).Value = "DOCUMENT_TEMPLATE_MANAGEMENT"
$ws.Cells.Item(110, 2).Value = "DOCUMENT_TEMPLATE_MANAGEMENT"

# Column C ("Admin") is the only "Yes" for this new right; everything else is "No".
$ws.Cells.Item(110, 3).Value = "Yes"
for ($col = 4; $col -le 25; $col++) {
    $ws.Cells.Item(110, $col).Value = "No"
}

# Bump the documented SORMAS version on the "About" sheet.
$ws2 = $wb.Worksheets.Item("About")
$ws2.Cells.Item(2, 1).Value = "1.50.0-SNAPSHOT"
